$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (MI): B2 formula changes from 2*D2 to D2*2, C2/D2 values change
$ws.Range("B2").Formula = "=D2*2"
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = 10

# Update row 3 (CSK): formula text changes (2*D3 -> D3*2), values change
$ws.Range("B3").Formula = "=D3*2"
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = 6

# Update row 4 (DC): formula text + values change
$ws.Range("B4").Formula = "=D4*2"
$ws.Range("C4").Value = 10
$ws.Range("D4").Value = 5

# Update row 5 (RCB): formula text + values change
$ws.Range("B5").Formula = "=D5*2"
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = 2

# Add new shared strings / rows 6-8: RR, KKR, SRH
$ws.Range("A6").Value = "RR"
$ws.Range("B6").Formula = "=D6*2"
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = 1
$ws.Range("E6").Formula = "=C6-D6"

$ws.Range("A7").Value = "KKR"
$ws.Range("B7").Formula = "=D7*2"
$ws.Range("C7").Value = 10
$ws.Range("D7").Value = 7
$ws.Range("E7").Formula = "=C7-D7"

$ws.Range("A8").Value = "SRH"
$ws.Range("B8").Formula = "=D8*2"
$ws.Range("C8").Value = 10
$ws.Range("D8").Value = 9
$ws.Range("E8").Formula = "=C8-D8"

# Update sheet view: zoom and selection
$ws.Application.ActiveWindow.Zoom = 256
$ws.Range("C9").Select()

$wb.Save()
